# split_main_means_nolegend.xlsx
# - Replace CRLF-wrapped "**test**" line-break markers in the shared-string
#   labels with literal "<br>" tags, all on one line (no more embedded
#   newline / **test** marker).
# - Small numeric tweaks to the CI_low / CI_high (columns C / D) values for
#   the "All" and "Europe" rows on every facet (re-run with a newer ggplot2 /
#   RNG on a different machine moved the bootstrap CIs very slightly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text labels (column E) -------------------------------------------------
# Every cell that shares a given label gets the same new text so the
# workbook keeps a single de-duplicated shared string per label.

$ws.Range("E2:E12").Value2  = "Global: Education, Healthcare and<br>Renewable energy in low-income countries"
$ws.Range("E24:E34").Value2 = "Share allocated to Global spending options<br>when 5 out of 13 options are randomly selected<br>(4 out of 13 being of Global nature)"
$ws.Range("E35:E45").Value2 = "Global: Education and Healthcare<br>in low-income countries"
$ws.Range("E46:E56").Value2 = "Global: Renewable energy and<br>infrastructure to cope with climate change"
$ws.Range("E57:E67").Value2 = "Global: Loss and Damage Fund (to<br>rebuild after climate disasters)"

# --- Numeric tweaks (columns C = CI_low, D = CI_high) -----------------------

$ws.Range("C2").Value2  = 17.1175986105978
$ws.Range("D2").Value2  = 17.8381270057554
$ws.Range("C3").Value2  = 18.1092138950899
$ws.Range("D3").Value2  = 19.2394096308851

$ws.Range("C13").Value2 = 25.4207830126572
$ws.Range("D13").Value2 = 26.2145778603093
$ws.Range("C14").Value2 = 28.0556145330742
$ws.Range("D14").Value2 = 29.2819146425341

$ws.Range("C24").Value2 = 26.2840599100719
$ws.Range("D24").Value2 = 27.4936456729721
$ws.Range("C25").Value2 = 26.9617128666953
$ws.Range("D25").Value2 = 28.7918085096094

$ws.Range("C35").Value2 = 17.9671758541937
$ws.Range("D35").Value2 = 19.2009014124105
$ws.Range("C36").Value2 = 18.6186128419952
$ws.Range("D36").Value2 = 20.4680649318187

$ws.Range("C46").Value2 = 17.7086767334779
$ws.Range("D46").Value2 = 18.8494959601504
$ws.Range("C47").Value2 = 17.8436710027609
$ws.Range("D47").Value2 = 19.5229413950389

$ws.Range("C57").Value2 = 16.1687738455109
$ws.Range("D57").Value2 = 17.2758576241267
$ws.Range("C58").Value2 = 15.6153160964197
$ws.Range("D58").Value2 = 17.2511550460299

$ws.Range("C68").Value2 = 15.7116001581493
$ws.Range("D68").Value2 = 16.7695833387339
$ws.Range("C69").Value2 = 16.4312908880532
$ws.Range("D69").Value2 = 17.9770834330384
